$wb = $excel.ActiveWorkbook

# Values in column F that changed, identical for both sheets except row 27
$updates = @{
    5  = 254
    6  = 44
    7  = 157
    8  = 256
    9  = 27
    13 = 28
    14 = 89
    15 = 758
    16 = 47
    17 = 485
    18 = 420
    19 = 139
    20 = 66
    21 = 33
    23 = 1270
    24 = 2924
    27 = 756
    28 = 64
    29 = 1635
    31 = 455
    32 = 23
    33 = 271
    34 = 395
    36 = 614
    37 = 424
    38 = 10
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
